$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header cell in column A changed from "User Story ID" to "ID"
$ws.Range("A1").Value = "ID"

# The active/selected cell on the sheet (bottom-right frozen pane) moved to A5
[void]$ws.Range("A5").Select()
